$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2 and 3
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 59

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 40

# Copy style from A3 to A4 and A5 (bold/bordered/centered style used in column A)
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Add new rows 4 and 5
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 37

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 28
